# Roster update: remove Marco Eigenmann and Dea Putri, add Laura-Florina
# Krattinger (Nestle), shrink the query table / named range from 16 to 15
# rows, and move the selection to N11:N12 (matches the author's last
# on-screen selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$lo = $ws.ListObjects.Item(1)

# --- remove "Dea Putri" (last data row, row 16: A16/B16) -------------------
$ws.Rows.Item(16).Delete()

# --- remove "Marco Eigenmann" (row 9: A9/B9) --------------------------------
$ws.Rows.Item(9).Delete()

# --- append the new committee member ---------------------------------------
$newRow = $lo.ListRows.Add()
$newCells = $newRow.Range
$newCells.Cells.Item(1,1).Value = "Laura-Florina"
$newCells.Cells.Item(1,2).Value = "Krattinger"
$newCells.Cells.Item(1,3).Value = "Nestlé"
$newCells.Cells.Item(1,4).NumberFormat = "General"
$newCells.Cells.Item(1,5).Value = 1

# --- fix up the hidden ExternalData_1 defined name (Sheet2!$A$1:$E$16 -> $E$15)
$n = $wb.Names.Item("ExternalData_1")
$n.RefersTo = "=Sheet2!`$A`$1:`$E`$15"

# --- leave a couple of formatted-but-empty cells the way the author did ----
$ws.Cells.Item(11,14).NumberFormat = "General"
$ws.Cells.Item(12,14).NumberFormat = "General"

# --- restore the on-screen selection ---------------------------------------
$ws.Range("N11:N12").Select()
